$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date the row was last refreshed.
# Bump it from 45175 (2023-09-06) to 45183 (2023-09-14) for every data
# row (rows 2 through 45), leaving the existing date formatting intact.
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
